# Auto-update draw results: append the 2025-09-24 Pick 3 draw as a new
# row at the bottom of the results table (mirrors the nightly results
# importer that appends one row per day).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet stores every column as literal text (dates, the zero-padded
# "phase" code and the timestamp are all text, not real numbers/dates).
# A leading apostrophe forces Excel to keep "2025-09-24" and "250924" as
# text instead of auto-converting them to a date serial / number.
$ws.Range("A8").Value = "'2025-09-24"
$ws.Range("B8").Value = "Pick 3"
$ws.Range("C8").Value = "'250924"
$ws.Range("D8").Value = "6-5-9"
$ws.Range("E8").Value = "2025-09-24T21:38:26.111+04:00"

# Typing a quote-prefixed value flags the cell with a "number stored as
# text" style; reset those two cells back to the plain/default style so
# the new row matches the rest of the (unstyled) table.
$ws.Range("A8").Style = "Normal"
$ws.Range("C8").Style = "Normal"
